# Revert "Merge pull request #48 from LakeFishing/main"
#
# The merge had replaced the "time" label used in a few cells with a
# placeholder value "special". This reverts those three cells back to
# "time" (matching the already-existing "time" labels in B6/B11 of the
# same sheet). Once nothing references the shared string "special"
# anymore, Excel drops the now-unused shared-string entries automatically
# when the workbook is re-saved, which is exactly what the target diff
# shows for xl/sharedStrings.xml (uniqueCount 34 -> 32, the two "special"
# <si> entries disappear).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- revert the three cell values that were changed to "special" --------
$ws.Range("F4").Value  = "time"
$ws.Range("B9").Value  = "time"
$ws.Range("B14").Value = "time"

# --- restore the saved selection (previously saved selection spanned
#     C7 and C12, with the active cell on C12; select C12 last so it is
#     the cell recorded as the active one) -------------------------------
$ws.Range("C7").Select()  | Out-Null
$ws.Range("C12").Select() | Out-Null
